# Proyectos.xlsx — Hoja1 fixes
#   1) Fix a duplicated catalogue value that only differed by a trailing
#      space ("CANALES ALTERNATIVOS " -> "CANALES ALTERNATIVOS") in I7.
#   2) Remove the blank/placeholder project row (old row 50) and let the
#      remaining rows shift up.
#   3) Re-apply AutoFilter over the (now smaller) table range and leave the
#      selection/scroll position where the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Activate()

# 1) Deduplicate "CANALES ALTERNATIVOS " (trailing space) -> "CANALES ALTERNATIVOS"
$ws.Range("I7").Value = "CANALES ALTERNATIVOS"

# 2) Delete the blank project row (row 50); rows below shift up automatically.
$ws.Rows.Item(50).Delete()

# 3) Rebuild the AutoFilter over the new used range (A1:L55) and register the
#    corresponding hidden _FilterDatabase name, mirroring what Excel does
#    when you press Ctrl+Shift+L on a selection inside the table.
$ws.Range("A1:L55").AutoFilter() | Out-Null
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=Hoja1!`$A`$1:`$L`$55")
$filterName.Visible = $false

# 4) Restore the saved view state: scrolled down so row 29 is at the top,
#    with M30 as the active/selected cell.
$ws.Range("M30").Select()
$excel.ActiveWindow.ScrollRow = 29
$excel.ActiveWindow.ScrollColumn = 1
